# Update the "Förändrad" (Changed) date column (column C) for all data rows.
# Every data row's value moves from serial date 45181 (2023-09-12) to
# 45182 (2023-09-13).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45181) {
        $cell.Value2 = 45182
    }
}
